# Apply the "add 2022-Q3 data" edit to the 688333 workbook.
$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1) "总计" sheet: insert a new row 2 for the 2022-Q3 summary figures.
#    Existing rows 2-6 (2022-Q2 .. 2020-Q4) shift down to rows 3-7.
# ------------------------------------------------------------------
$wsTotal = $wb.Worksheets.Item(1)
$wsTotal.Rows.Item(2).Insert()

# Clear any formatting the row-insert may have copied down so the new
# row matches the plain look of the other data rows.
$wsTotal.Range("A2:D2").ClearFormats()

# Re-apply the same style used by the other "序号" cells in column A.
$wsTotal.Range("A3").Copy()
$wsTotal.Range("A2").PasteSpecial(-4122)

$wsTotal.Cells.Item(2, 1).Value = 0
$wsTotal.Cells.Item(2, 2).Value = "2022-Q3"
$wsTotal.Cells.Item(2, 3).Value = 19
$wsTotal.Cells.Item(2, 4).Value = 6.06

# The "A" column is a plain 0-based row index, not data that moves with
# the row - renumber the rows that shifted down (rows 3-7).
$wsTotal.Cells.Item(3, 1).Value = 1
$wsTotal.Cells.Item(4, 1).Value = 2
$wsTotal.Cells.Item(5, 1).Value = 3
$wsTotal.Cells.Item(6, 1).Value = 4
$wsTotal.Cells.Item(7, 1).Value = 5

# ------------------------------------------------------------------
# 2) Insert a brand-new "2022-Q3" sheet right after "总计" (i.e. right
#    before the existing "2022-Q2" sheet), holding the fund holdings
#    detail for the new quarter.
# ------------------------------------------------------------------
$wsBefore = $wb.Worksheets.Item(2)
$ws = $wb.Worksheets.Add($wsBefore)
$ws.Name = "2022-Q3"

$headers = @("基金代码","基金名称","基金规模","股票总仓位","仓位占比","持有市值(亿元)","仓位排名")
for ($c = 0; $c -lt $headers.Length; $c++) {
    $ws.Cells.Item(1, $c + 2).Value = $headers[$c]
}
# Header row uses the same bold/bordered style as the other sheets.
$wsTotal.Range("B1").Copy()
$ws.Range("B1:H1").PasteSpecial(-4122)
for ($c = 0; $c -lt $headers.Length; $c++) {
    $ws.Cells.Item(1, $c + 2).Value = $headers[$c]
}

$rows = @(
    @("501207","华夏创新未来混合（LOF）","44.72","91.75","4.07","1.8201",10),
    @("000031","华夏复兴混合A","24.25","88.23","4.32","1.0476",9),
    @("160106","南方高增长混合（LOF）","15.29","91.51","5.65","0.8639",3),
    @("007349","华夏科技创新混合A","11.68","90.65","4.20","0.4906",7),
    @("160105","南方积极配置混合（LOF）","5.42","91.29","5.65","0.3062",4),
    @("013962","华夏创新视野一年持有混合A","7.47","88.70","3.99","0.2981",10),
    @("005358","东方阿尔法精选灵活配置混合A","3.47","93.91","8.17","0.2835",6),
    @("013963","华夏创新视野一年持有混合C","6.46","88.70","3.99","0.2578",10),
    @("010518","华夏先锋科技一年定期开放混合A","5.64","89.54","4.20","0.2369",8),
    @("001173","中欧瑾和灵活配置混合 - A","2.26","92.00","4.65","0.1051",7),
    @("000554","南方中国梦灵活配置混合","1.33","90.51","5.89","0.0783",3),
    @("010519","华夏先锋科技一年定期开放混合C","1.78","89.54","4.20","0.0748",8),
    @("002577","南方新兴龙头灵活配置混合","1.52","78.68","4.68","0.0711",6),
    @("005359","东方阿尔法精选灵活配置混合C","0.58","93.91","8.17","0.0474",6),
    @("007350","华夏科技创新混合C","0.88","90.65","4.20","0.0370",7),
    @("012669","南方新兴产业混合A","0.53","81.22","3.49","0.0185",8),
    @("012670","南方新兴产业混合C","0.40","81.22","3.49","0.0140",8),
    @("001174","中欧瑾和灵活配置混合 - C","0.23","92.00","4.65","0.0107",7),
    @("015073","华夏复兴混合C","0.04","88.23","4.32","0.0017",9)
)

$r = 2
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = $r - 2
    $wsTotal.Range("A2").Copy()
    $ws.Cells.Item($r, 1).PasteSpecial(-4122)
    $ws.Cells.Item($r, 1).Value = $r - 2

    $ws.Cells.Item($r, 2).Value = "'" + $row[0]
    $ws.Cells.Item($r, 3).Value = "'" + $row[1]
    $ws.Cells.Item($r, 4).Value = "'" + $row[2]
    $ws.Cells.Item($r, 5).Value = "'" + $row[3]
    $ws.Cells.Item($r, 6).Value = "'" + $row[4]
    $ws.Cells.Item($r, 7).Value = "'" + $row[5]
    $ws.Range($ws.Cells.Item($r, 2), $ws.Cells.Item($r, 7)).ClearFormats()

    $ws.Cells.Item($r, 8).Value = $row[6]

    $r = $r + 1
}

Write-Host "2022-Q3 sheet and zongji summary row added"
